$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value2 = "ORD-1736742537979-61cd6af4"
$ws.Range("B8").Value2 = "Prithviraj"
$ws.Range("C8").Value2 = "iamprithvi@gmail.com"
$ws.Range("D8").Value2 = 699
$ws.Range("E8").Value2 = 0
$ws.Range("F8").Value2 = "null"
$ws.Range("G8").Value2 = 699
$ws.Range("H8").Value2 = "pending"
$ws.Range("I8").Value2 = "13/01/2025, 09:58:57"

# Row 9
$ws.Range("A9").Value2 = "ORD-1736742796596-a5238dcc"
$ws.Range("B9").Value2 = "Prithviraj"
$ws.Range("C9").Value2 = "iamprithvi@gmail.com"
$ws.Range("D9").Value2 = 1518
$ws.Range("E9").Value2 = 0
$ws.Range("F9").Value2 = "null"
$ws.Range("G9").Value2 = 1518
$ws.Range("H9").Value2 = "pending"
$ws.Range("I9").Value2 = "13/01/2025, 10:03:16"

# Row 10
$ws.Range("A10").Value2 = "ORD-1736742919832-eddab4ba"
$ws.Range("B10").Value2 = "Prithviraj"
$ws.Range("C10").Value2 = "iamprithvi@gmail.com"
$ws.Range("D10").Value2 = 759
$ws.Range("E10").Value2 = 0
$ws.Range("F10").Value2 = "null"
$ws.Range("G10").Value2 = 759
$ws.Range("H10").Value2 = "pending"
$ws.Range("I10").Value2 = "13/01/2025, 10:05:19"

# Row 11
$ws.Range("A11").Value2 = "ORD-1736743627086-35daa2b2"
$ws.Range("B11").Value2 = "Prithviraj"
$ws.Range("C11").Value2 = "iamprithvi@gmail.com"
$ws.Range("D11").Value2 = 709
$ws.Range("E11").Value2 = 128
$ws.Range("F11").Value2 = "MAX50"
$ws.Range("G11").Value2 = 581
$ws.Range("H11").Value2 = "canceled"
$ws.Range("I11").Value2 = "13/01/2025, 10:17:07"

# Row 12
$ws.Range("A12").Value2 = "ORD-1736743888962-efc6cfe5"
$ws.Range("B12").Value2 = "Prithviraj"
$ws.Range("C12").Value2 = "iamprithvi@gmail.com"
$ws.Range("D12").Value2 = 759
$ws.Range("E12").Value2 = 0
$ws.Range("F12").Value2 = "null"
$ws.Range("G12").Value2 = 759
$ws.Range("H12").Value2 = "canceled"
$ws.Range("I12").Value2 = "13/01/2025, 10:21:28"

# Row 13
$ws.Range("A13").Value2 = "ORD-1736744544825-5382ce92"
$ws.Range("B13").Value2 = "Prithviraj"
$ws.Range("C13").Value2 = "iamprithvi@gmail.com"
$ws.Range("D13").Value2 = 711
$ws.Range("E13").Value2 = 0
$ws.Range("F13").Value2 = "null"
$ws.Range("G13").Value2 = 711
$ws.Range("H13").Value2 = "pending"
$ws.Range("I13").Value2 = "13/01/2025, 10:32:24"
